$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: duplicate of the Dubai (DSC) / Sunrisers match entry
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " October 27 2020"
$ws.Range("C4").Value = "Sunrisers won by 88 runs"
$ws.Range("D4").Value = "Delhi Capitals"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
$ws.Range("F4").Value = "Tushar Deshpande "

$ws.Range("G4:K4").NumberFormat = "@"
$ws.Range("G4").Value = "20"
$ws.Range("H4").Value = "9"
$ws.Range("I4").Value = "2"
$ws.Range("J4").Value = "1"
$ws.Range("K4").Value = "222.22"
$ws.Range("G4:K4").Style = "Normal"

# Row 5: duplicate of the Abu Dhabi / KKR match entry
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " October 24 2020"
$ws.Range("C5").Value = "KKR won by 59 runs"
$ws.Range("D5").Value = "Delhi Capitals"
$ws.Range("E5").Value = "Kolkata Knight Riders"
$ws.Range("F5").Value = "Tushar Deshpande "

$ws.Range("G5:K5").NumberFormat = "@"
$ws.Range("G5").Value = "1"
$ws.Range("H5").Value = "3"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "33.33"
$ws.Range("G5:K5").Style = "Normal"
